$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'95.374.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.30%  '

# Row 3
$ws.Range("D3").Value = "'3.616.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.87%  '

# Row 4
$ws.Range("D4").Value = "'2.73"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +26.08%  '

# Row 5
$ws.Range("E5").Value = '  +0.03%  '

# Row 6
$ws.Range("D6").Value = "'224.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.85%  '

# Row 7
$ws.Range("D7").Value = "'641.85"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.37%  '

# Row 8
$ws.Range("D8").Value = "'0.422"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.71%  '

# Row 9
$ws.Range("E9").Value = '  +5.51%  '

# Row 10
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

# Row 11
$ws.Range("D11").Value = "'3.614.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.89%  '

# Row 12
$ws.Range("D12").Value = "'50.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.39%  '

# Row 13
$ws.Range("D13").Value = "'0.218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.04%  '

# Row 14
$ws.Range("E14").Value = '  -6.40%  '

# Row 15
$ws.Range("E15").Value = '  -5.13%  '

# Row 16
$ws.Range("D16").Value = "'4.289.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.06%  '

# Row 17
$ws.Range("D17").Value = "'95.184.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.38%  '

# Row 18
$ws.Range("D18").Value = "'24.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +28.90%  '

# Row 19
$ws.Range("D19").Value = "'9.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '

# Row 20
$ws.Range("E20").Value = '  +4.66%  '

# Row 21
$ws.Range("D21").Value = "'3.614.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.01%  '

# Row 22
$ws.Range("D22").Value = "'0.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +37.86%  '

# Row 23
$ws.Range("D23").Value = "'0.536"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.61%  '

# Row 24
$ws.Range("D24").Value = "'137.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +15.61%  '

# Row 25
$ws.Range("D25").Value = "'532.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
$ws.Range("E26").Value = '  -5.49%  '

# Row 27
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = "'0.0000203"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.42%  '

# Row 28
$ws.Range("B28").Value = 'NEARProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D28").Value = "'7.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.28%  '

# Row 29
$ws.Range("D29").Value = "'13.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.17%  '

# Row 30
$ws.Range("D30").Value = "'3.786.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.51%  '

# Row 31
$ws.Range("D31").Value = "'13.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.46%  '

# Row 32
$ws.Range("D32").Value = "'3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.79%  '

# Row 33
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("D34").Value = "'0.645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.40%  '

# Row 35
$ws.Range("D35").Value = "'1.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.47%  '

# Row 36
$ws.Range("D36").Value = "'33.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.29%  '

# Row 37
$ws.Range("D37").Value = "'0.183"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.88%  '

# Row 38
$ws.Range("E38").Value = '  +0.25%  '

# Row 39
$ws.Range("D39").Value = "'0.0563"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +22.27%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = "'7.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.39%  '

# Row 42
$ws.Range("D42").Value = "'8.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.80%  '

# Row 43
$ws.Range("D43").Value = "'593.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.54%  '

# Row 44
$ws.Range("D44").Value = "'0.513"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.56%  '

# Row 45
$ws.Range("D45").Value = "'1.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.81%  '

# Row 46
$ws.Range("D46").Value = "'40.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.46%  '

# Row 47
$ws.Range("E47").Value = '  +0.02%  '

# Row 48
$ws.Range("E48").Value = '  -7.20%  '

# Row 49
$ws.Range("D49").Value = "'9.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.99%  '

# Row 50
$ws.Range("D50").Value = "'236.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.82%  '

# Row 51
$ws.Range("D51").Value = "'2.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.89%  '
